# Update the handback status report with freshly generated timestamps.
$wb = $excel.ActiveWorkbook

# zh-cn sheet: row for 879e6e22-...zh-cn.xlf entry (row 3)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-01-18 02:53:22"
$wsZhCn.Range("G3").Value = "2016-01-18 02:54:03"

# de-de sheet: row for 879e6e22-...de-de.xlf entry (row 3)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-01-18 02:53:33"
$wsDeDe.Range("G3").Value = "2016-01-18 02:54:21"
